$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value as text without altering its cell style,
# even when the string would otherwise be auto-parsed as a number by Excel.
function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "67.201.60"
Set-TextValue $ws.Range("E2") "  +0.11%  "

# Row 3
Set-TextValue $ws.Range("D3") "3.132.46"
Set-TextValue $ws.Range("E3") "  +0.35%  "

# Row 5
Set-TextValue $ws.Range("E5") "  -0.10%  "

# Row 6
Set-TextValue $ws.Range("D6") "177.83"
Set-TextValue $ws.Range("E6") "  +1.53%  "

# Row 7
Set-TextValue $ws.Range("E7") "  -0.09%  "

# Row 8
Set-TextValue $ws.Range("D8") "3.136.67"
Set-TextValue $ws.Range("E8") "  +0.62%  "

# Row 9
Set-TextValue $ws.Range("D9") "0.518"
Set-TextValue $ws.Range("E9") "  -0.91%  "

# Row 10
Set-TextValue $ws.Range("D10") "6.43"
Set-TextValue $ws.Range("E10") "  +0.18%  "

# Row 11
Set-TextValue $ws.Range("D11") "0.153"
Set-TextValue $ws.Range("E11") "  -1.61%  "

# Row 12
Set-TextValue $ws.Range("D12") "0.477"
Set-TextValue $ws.Range("E12") "  -0.72%  "

# Row 13
Set-TextValue $ws.Range("D13") "0.0000243"
Set-TextValue $ws.Range("E13") "  -2.40%  "

# Row 14
Set-TextValue $ws.Range("D14") "36.57"
Set-TextValue $ws.Range("E14") "  -1.74%  "

# Row 15
Set-TextValue $ws.Range("E15") "  -0.29%  "

# Row 16
Set-TextValue $ws.Range("D16") "3.651.37"
Set-TextValue $ws.Range("E16") "  +0.33%  "

# Row 17
Set-TextValue $ws.Range("D17") "67.161.59"
Set-TextValue $ws.Range("E17") "  +0.06%  "

# Row 18
Set-TextValue $ws.Range("B18") "Polkadot"
Set-TextValue $ws.Range("C18") "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws.Range("D18") "7.06"
Set-TextValue $ws.Range("E18") "  -0.94%  "

# Row 19
Set-TextValue $ws.Range("B19") "Chainlink"
Set-TextValue $ws.Range("C19") "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue $ws.Range("D19") "16.92"
Set-TextValue $ws.Range("E19") "  +2.81%  "

# Row 20
Set-TextValue $ws.Range("D20") "3.129.90"
Set-TextValue $ws.Range("E20") "  +0.26%  "

# Row 21
Set-TextValue $ws.Range("D21") "489.13"
Set-TextValue $ws.Range("E21") "  -0.42%  "

# Row 22
Set-TextValue $ws.Range("E22") "  -1.07%  "

# Row 23
Set-TextValue $ws.Range("D23") "0.697"
Set-TextValue $ws.Range("E23") "  -1.48%  "

# Row 24
Set-TextValue $ws.Range("D24") "84.13"
Set-TextValue $ws.Range("E24") "  -0.11%  "

# Row 25
Set-TextValue $ws.Range("D25") "12.92"
Set-TextValue $ws.Range("E25") "  -2.15%  "

# Row 26
Set-TextValue $ws.Range("D26") "2.27"
Set-TextValue $ws.Range("E26") "  -1.55%  "

# Row 27
Set-TextValue $ws.Range("D27") "10.28"
Set-TextValue $ws.Range("E27") "  -1.00%  "

# Row 28
Set-TextValue $ws.Range("E28") "  +0.08%  "

# Row 29
Set-TextValue $ws.Range("D29") "8.09"
Set-TextValue $ws.Range("E29") "  +2.18%  "

# Row 30
Set-TextValue $ws.Range("D30") "2.31"
Set-TextValue $ws.Range("E30") "  -1.93%  "

# Row 31
Set-TextValue $ws.Range("E31") "  -2.69%  "

# Row 32
Set-TextValue $ws.Range("D32") "28.27"
Set-TextValue $ws.Range("E32") "  -1.24%  "

# Row 33
Set-TextValue $ws.Range("D33") "0.113"
Set-TextValue $ws.Range("E33") "  -1.17%  "

# Row 34
Set-TextValue $ws.Range("D34") "0.0₃0947"
Set-TextValue $ws.Range("E34") "  -0.38%  "

# Row 36
Set-TextValue $ws.Range("D36") "48.83"
Set-TextValue $ws.Range("E36") "  +3.56%  "

# Row 37
Set-TextValue $ws.Range("D37") "5.66"
Set-TextValue $ws.Range("E37") "  -3.78%  "

# Row 38
Set-TextValue $ws.Range("D38") "0.949"
Set-TextValue $ws.Range("E38") "  -2.75%  "

# Row 39
Set-TextValue $ws.Range("E39") "  +0.59%  "

# Row 40
Set-TextValue $ws.Range("D40") "49.49"
Set-TextValue $ws.Range("E40") "  -1.26%  "

# Row 41
Set-TextValue $ws.Range("B41") "Stacks"
Set-TextValue $ws.Range("C41") "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws.Range("D41") "2.01"
Set-TextValue $ws.Range("E41") "  -2.47%  "

# Row 42
Set-TextValue $ws.Range("B42") "Kaspa"
Set-TextValue $ws.Range("C42") "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws.Range("D42") "0.124"
Set-TextValue $ws.Range("E42") "  -0.11%  "

# Row 43
Set-TextValue $ws.Range("D43") "8.38"
Set-TextValue $ws.Range("E43") "  -1.68%  "

# Row 44
Set-TextValue $ws.Range("D44") "2.71"
Set-TextValue $ws.Range("E44") "  +4.07%  "

# Row 45
Set-TextValue $ws.Range("D45") "2.813.66"
Set-TextValue $ws.Range("E45") "  -0.23%  "

# Row 46
Set-TextValue $ws.Range("D46") "377.53"
Set-TextValue $ws.Range("E46") "  -1.65%  "

# Row 47
Set-TextValue $ws.Range("D47") "0.0351"
Set-TextValue $ws.Range("E47") "  -0.96%  "

# Row 48
Set-TextValue $ws.Range("D48") "134.83"
Set-TextValue $ws.Range("E48") "  -0.52%  "

# Row 50
Set-TextValue $ws.Range("D50") "24.93"
Set-TextValue $ws.Range("E50") "  -0.05%  "

# Row 51
Set-TextValue $ws.Range("D51") "2.25"
Set-TextValue $ws.Range("E51") "  +1.56%  "
